$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the totals row (row 38) for the new product line,
# then copy the formatting (incl. borders/merges) of the previous item row
# (row 37) into it so it matches the existing table pattern exactly.
$ws.Rows("38").Insert()
$ws.Range("A37:N37").Copy($ws.Range("A38:N38"))
$excel.CutCopyMode = $false
$ws.Rows("38").RowHeight = 25.5

# New item values (item #35)
$ws.Range("A38").Value = 35
$ws.Range("B38").Value = "معطر فريدا "
$ws.Range("H38").Value = "9:0"
$ws.Range("L38").Value = 65
# N38 keeps the same category text as N37 ("بلاستر مترسيلك 2.5 سم"), already copied.

# Update the totals row (now shifted to row 39): add new item's value to the total.
$ws.Range("K39").Value = 2118.6399999999999

# Footer row (now shifted to row 40) height changed slightly.
$ws.Rows("40").RowHeight = 16.5
